$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.075.80"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.649.87"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.94"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5273"
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2608"
$ws.Range("E8").Value = "  -2.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06305"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.34"
$ws.Range("E10").Value = "  -3.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07749"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.468"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Value = "1.623.49"
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5449"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").Value = "0.0₅8100"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.12"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "26.099.74"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.542"
$ws.Range("E19").Value = "  -2.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.83"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.03"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.977"
$ws.Range("E22").Value = "  -1.89%  "
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "139.93"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1238"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.235"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.16"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.434"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05902"
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.278"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("E31").Value = "  -2.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.233"
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("E33").Value = "  -6.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.413"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9423"
$ws.Range("E35").Value = "  -4.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.755"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5655"
$ws.Range("E37").Value = "  -4.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01603"
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.837"
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8451"
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.003"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.70"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").Value = "1.006.42"
$ws.Range("E43").Value = "  -3.28%  "
$ws.Range("D44").Value = "1.797.92"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "56.82"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").Value = "0.0₈106"
$ws.Range("E46").Value = "  -4.07%  "
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4289"
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.476"
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05151"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.791"
$ws.Range("E51").Value = "  -3.76%  "
